# "You can't place roads on towers anymore!"
# Add a new "Speed" stat column (F) to the tower comparison table (rows 34-45)
# and bump several tower Health (column D) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column F
$ws.Range("F34").Value = "Speed"

# Row 35 - Health 35 -> 50, new Speed = 1
$ws.Range("D35").Value = 50
$ws.Range("F35").Value = 1

# Row 36 - Health 50 -> 100, new Speed = 1
$ws.Range("D36").Value = 100
$ws.Range("F36").Value = 1

# Row 37 - Health 20 -> 1000, new Speed = 1
$ws.Range("D37").Value = 1000
$ws.Range("F37").Value = 1

# Row 38 - Health 25 -> 100, new Speed = 2
$ws.Range("D38").Value = 100
$ws.Range("F38").Value = 2

# Row 39 - Health 60 -> 80, new Speed = 1
$ws.Range("D39").Value = 80
$ws.Range("F39").Value = 1

# Row 40 - Health 35 -> 60, new Speed = 1
$ws.Range("D40").Value = 60
$ws.Range("F40").Value = 1

# Row 41 - Health 20 -> 50, new Speed = 3
$ws.Range("D41").Value = 50
$ws.Range("F41").Value = 3

# Row 42 - Health 35 -> 50, new Speed = 2
$ws.Range("D42").Value = 50
$ws.Range("F42").Value = 2

# Row 43 - Health 45 -> 150, new Speed = 1
$ws.Range("D43").Value = 150
$ws.Range("F43").Value = 1

# Row 44 - Health 35 -> 80, new Speed = 1
$ws.Range("D44").Value = 80
$ws.Range("F44").Value = 1

# Row 45 - Health 30 -> 40, new Speed = 1
$ws.Range("D45").Value = 40
$ws.Range("F45").Value = 1

# Update the view: scroll so row 10 area is visible and select F46
$ws.Range("F46").Select()
